$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" (empty value) right after "Contact" (row 10)
# and before "Description" (row 11), pushing Description and subsequent rows down.
$ws.Rows("11").Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
